$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "C6, C10"
$ws.Range("B18").Value = "Decoupling Capacitor"

$ws.Range("B19").Value = "Decoupling Capacitor"

$ws.Range("B20").Value = "Decoupling Capacitor"

$ws.Range("A21").Value = "S1"

$ws.Range("A22").Value = "R7"
$ws.Range("B22").Value = "HV Divider"

$ws.Range("A23").Value = "R9"
$ws.Range("B23").Value = "HV Divider"

$ws.Range("A24").Value = "R8"
$ws.Range("B24").Value = "HV Divider"

$ws.Range("C17:E24").ClearFormats() | Out-Null

$ws.Range("H25").Select() | Out-Null
